$wb = $excel.ActiveWorkbook

function Update-AnimeSheet($ws, [int]$lastRow) {

    # --- numeric "want to go" count bumps on existing rows (before the insert) ---
    $ws.Range("F2").Value = 230
    $ws.Range("F5").Value = 811
    $ws.Range("F7").Value = 6417
    $ws.Range("F8").Value = 49
    $ws.Range("F9").Value = 70
    $ws.Range("F10").Value = 109
    $ws.Range("F11").Value = 70

    # --- insert a new row 14 (new event: 2024-05-18 偶活企划) ---
    $ws.Rows.Item(14).Insert()

    # Copy formats from the row above (13) onto the freshly inserted row so
    # the bold/centered/bordered style used by column A (and the rest of
    # the row) matches the sheet's existing look.
    $ws.Range("A13:I13").Copy()
    $ws.Range("A14:I14").PasteSpecial(-4122)
    $excel.CutCopyMode = $false

    $ws.Range("A14").Value = 13
    $ws.Range("B14").NumberFormat = "@"
    $ws.Range("B14").Value = "2024-05-18"
    $ws.Range("C14").Value = "合肥·首届偶活企划——偶像计划-闪耀舞台"
    $ws.Range("D14").Value = "阜阳路16号 银瑞林国际大酒店"
    $ws.Range("E14").Value = "2024.05.18 09:00-05.18 17:00"
    $ws.Range("F14").Value = 1
    $ws.Range("G14").Value = 58
    $ws.Range("H14").Value = "https://show.bilibili.com/platform/detail.html?id=83891"
    $ws.Range("I14").Value = "//i2.hdslb.com/bfs/openplatform/202404/lfqv8l9Q1712453982625.jpeg"

    # --- rows that used to be 14..lastRow are now 15..(lastRow+1); their
    # running index in column A needs to be bumped by one to stay in sync ---
    $r = 15
    while ($r -le ($lastRow + 1)) {
        $idx = $ws.Range("A$r").Value2
        $ws.Range("A$r").Value = $idx + 1
        $r = $r + 1
    }

    # the old last row ("MAX特摄only展" in 展览 / the music event row in 全部类型)
    # also got its own "想去人数" bump on top of the shift
    $newLastRow = $lastRow + 1
    $curF = $ws.Range("F$newLastRow").Value2
    $ws.Range("F$newLastRow").Value = $curF + 1
}

# Sheet "展览" (exhibitions) — rows 1-16 -> 1-17
$ws1 = $wb.Worksheets.Item("展览")
Update-AnimeSheet $ws1 16

# Sheet "演出" (performances) — single data row bump
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 5

# Sheet "本地生活" (local life) — untouched

# Sheet "全部类型" (all types) — rows 1-17 -> 1-18
$ws4 = $wb.Worksheets.Item("全部类型")
Update-AnimeSheet $ws4 17

$excel.CutCopyMode = $false
